# Adds a "2022-Q4" sheet (with fund holding data) ahead of the existing
# "2022-Q3" sheet, and updates the "总计" (totals) summary sheet with a new
# leading row for 2022-Q4, shifting the previous rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: force a cell to hold a literal text value (not auto-coerced to a
# number), while leaving the cell's style at the workbook default ("Normal")
# once the value has been written.
# ---------------------------------------------------------------------
function Set-TextCell($cell, [string]$val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Helper: apply the same "header" look already used on the existing
# quarter sheets (bold font, thin box border, centered/top aligned).
# ---------------------------------------------------------------------
function Set-HeaderStyle($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $rng.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $rng.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $rng.Borders.Item(10).LineStyle = 1  # xlEdgeRight
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet immediately before the current
#    "2022-Q3" sheet (tab order becomes 总计, 2022-Q4, 2022-Q3, 2022-Q2).
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Add($q3Sheet)
$q4Sheet.Name = "2022-Q4"

# Header row
$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"
Set-HeaderStyle $q4Sheet.Range("B1:H1")

# Data rows: index column A, fund code, fund name, fund size, stock
# position, position ratio, held market value, position rank.
$rows = @(
    @(0, "161611", "融通内需驱动混合A/B",       "9.28", "92.61", "4.85", "0.4501", 6),
    @(1, "014109", "融通内需驱动混合C",          "3.63", "92.61", "4.85", "0.1761", 6),
    @(2, "001319", "农银汇理信息传媒主题股票",   "2.92", "84.93", "4.87", "0.1422", 9),
    @(3, "014106", "融通成长30灵活配置混合C",    "2.19", "94.02", "4.84", "0.1060", 6),
    @(4, "002252", "融通成长30灵活配置混合A/B",  "1.65", "94.02", "4.84", "0.0799", 6),
    @(5, "001223", "鹏华文化传媒娱乐股票",       "0.77", "83.57", "4.90", "0.0377", 4)
)

$r = 2
foreach ($row in $rows) {
    $aCell = $q4Sheet.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.Item(7).LineStyle = 1
    $aCell.Borders.Item(8).LineStyle = 1
    $aCell.Borders.Item(9).LineStyle = 1
    $aCell.Borders.Item(10).LineStyle = 1

    Set-TextCell $q4Sheet.Cells.Item($r, 2) $row[1]
    Set-TextCell $q4Sheet.Cells.Item($r, 3) $row[2]
    Set-TextCell $q4Sheet.Cells.Item($r, 4) $row[3]
    Set-TextCell $q4Sheet.Cells.Item($r, 5) $row[4]
    Set-TextCell $q4Sheet.Cells.Item($r, 6) $row[5]
    Set-TextCell $q4Sheet.Cells.Item($r, 7) $row[6]
    $q4Sheet.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row for 2022-Q4 above
#    the existing 2022-Q3 row, pushing the rest down by one.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 0.99

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2

Write-Host "done"
